$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the three numeric-looking "확정공모가" text cells as Text before writing,
# matching the stored string type (t="s") used for the non-numeric "-" cells in this column.
$ws.Range("D19:D21").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "에이직랜드"
$ws.Cells.Item(2, 2).Value = "2023.10.23~10.27"
$ws.Cells.Item(2, 3).Value = "19,100~21,400"
$ws.Cells.Item(2, 4).Value = "-"
$ws.Cells.Item(2, 5).Value = 50353
$ws.Cells.Item(2, 6).Value = "삼성증권"

$ws.Cells.Item(3, 1).Value = "비아이매트릭스"
$ws.Cells.Item(3, 2).Value = "2023.10.19~10.25"
$ws.Cells.Item(3, 3).Value = "9,100~11,000"
$ws.Cells.Item(3, 4).Value = "-"
$ws.Cells.Item(3, 5).Value = 10920
$ws.Cells.Item(3, 6).Value = "IBK투자증권"

$ws.Cells.Item(4, 1).Value = "유투바이오"
$ws.Cells.Item(4, 2).Value = "2023.10.18~10.19"
$ws.Cells.Item(4, 3).Value = "3,300~3,900"
$ws.Cells.Item(4, 4).Value = "-"
$ws.Cells.Item(4, 5).Value = 3724
$ws.Cells.Item(4, 6).Value = "신한투자증권"

$ws.Cells.Item(5, 1).Value = "큐로셀"
$ws.Cells.Item(5, 2).Value = "2023.10.18~10.24"
$ws.Cells.Item(5, 3).Value = "29,800~33,500"
$ws.Cells.Item(5, 4).Value = "-"
$ws.Cells.Item(5, 5).Value = 47680
$ws.Cells.Item(5, 6).Value = "미래에셋증권,삼성증권"

$ws.Cells.Item(6, 1).Value = "서울보증보험"
$ws.Cells.Item(6, 2).Value = "2023.10.13~10.19"
$ws.Cells.Item(6, 3).Value = "39,500~51,800"
$ws.Cells.Item(6, 4).Value = "-"
$ws.Cells.Item(6, 5).Value = 275795
$ws.Cells.Item(6, 6).Value = "미래에셋증권,삼성증권"

$ws.Cells.Item(7, 1).Value = "유진테크놀로지"
$ws.Cells.Item(7, 2).Value = "2023.10.11~10.17"
$ws.Cells.Item(7, 3).Value = "12,800~14,500"
$ws.Cells.Item(7, 4).Value = "-"
$ws.Cells.Item(7, 5).Value = 13433
$ws.Cells.Item(7, 6).Value = "NH투자증권"

$ws.Cells.Item(8, 1).Value = "퀄리타스반도체"
$ws.Cells.Item(8, 2).Value = "2023.10.06~10.13"
$ws.Cells.Item(8, 3).Value = "13,000~15,000"
$ws.Cells.Item(8, 4).Value = "-"
$ws.Cells.Item(8, 5).Value = 23400
$ws.Cells.Item(8, 6).Value = "한국투자증권"

$ws.Cells.Item(9, 1).Value = "컨텍"
$ws.Cells.Item(9, 2).Value = "2023.10.06~10.13"
$ws.Cells.Item(9, 3).Value = "20,300~22,500"
$ws.Cells.Item(9, 4).Value = "-"
$ws.Cells.Item(9, 5).Value = 41818
$ws.Cells.Item(9, 6).Value = "대신증권"

$ws.Cells.Item(10, 1).Value = "워트"
$ws.Cells.Item(10, 2).Value = "2023.10.05~10.12"
$ws.Cells.Item(10, 3).Value = "5,000~5,600"
$ws.Cells.Item(10, 4).Value = "-"
$ws.Cells.Item(10, 5).Value = 20000
$ws.Cells.Item(10, 6).Value = "키움증권"

$ws.Cells.Item(11, 1).Value = "신성에스티"
$ws.Cells.Item(11, 2).Value = "2023.09.22~10.04"
$ws.Cells.Item(11, 3).Value = "22,000~25,000"
$ws.Cells.Item(11, 4).Value = "-"
$ws.Cells.Item(11, 5).Value = 44000
$ws.Cells.Item(11, 6).Value = "미래에셋증권"

$ws.Cells.Item(12, 1).Value = "퓨릿(구.신디프)"
$ws.Cells.Item(12, 2).Value = "2023.09.20~09.26"
$ws.Cells.Item(12, 3).Value = "8,800~10,700"
$ws.Cells.Item(12, 4).Value = "-"
$ws.Cells.Item(12, 5).Value = 36405
$ws.Cells.Item(12, 6).Value = "미래에셋증권"

$ws.Cells.Item(13, 1).Value = "에이치엠씨아이비스팩6호"
$ws.Cells.Item(13, 2).Value = "2023.09.19~09.20"
$ws.Cells.Item(13, 3).Value = "2,000~2,000"
$ws.Cells.Item(13, 4).Value = "-"
$ws.Cells.Item(13, 5).Value = 8000
$ws.Cells.Item(13, 6).Value = "현대차증권"

$ws.Cells.Item(14, 1).Value = "에스엘에스바이오"
$ws.Cells.Item(14, 2).Value = "2023.09.18~09.22"
$ws.Cells.Item(14, 3).Value = "8,200~9,400"
$ws.Cells.Item(14, 4).Value = "-"
$ws.Cells.Item(14, 5).Value = 6314
$ws.Cells.Item(14, 6).Value = "하나증권"

$ws.Cells.Item(15, 1).Value = "신한스팩11호"
$ws.Cells.Item(15, 2).Value = "2023.09.14~09.15"
$ws.Cells.Item(15, 3).Value = "2,000~2,000"
$ws.Cells.Item(15, 4).Value = "-"
$ws.Cells.Item(15, 5).Value = 36000
$ws.Cells.Item(15, 6).Value = "신한투자증권"

$ws.Cells.Item(16, 1).Value = "레뷰코퍼레이션"
$ws.Cells.Item(16, 2).Value = "2023.09.11~09.15"
$ws.Cells.Item(16, 3).Value = "11,500~13,200"
$ws.Cells.Item(16, 4).Value = "-"
$ws.Cells.Item(16, 5).Value = 25760
$ws.Cells.Item(16, 6).Value = "삼성증권"

$ws.Cells.Item(17, 1).Value = "두산로보틱스"
$ws.Cells.Item(17, 2).Value = "2023.09.11~09.15"
$ws.Cells.Item(17, 3).Value = "21,000~26,000"
$ws.Cells.Item(17, 4).Value = "-"
$ws.Cells.Item(17, 5).Value = 340200
$ws.Cells.Item(17, 6).Value = "한국투자증권,미래에셋증권,NH투자증권,KB증권,크레디트스위스증권"

$ws.Cells.Item(18, 1).Value = "한싹"
$ws.Cells.Item(18, 2).Value = "2023.09.08~09.14"
$ws.Cells.Item(18, 3).Value = "8,900~11,000"
$ws.Cells.Item(18, 4).Value = "-"
$ws.Cells.Item(18, 5).Value = 13350
$ws.Cells.Item(18, 6).Value = "케이비증권"

$ws.Cells.Item(19, 1).Value = "밀리의서재"
$ws.Cells.Item(19, 2).Value = "2023.09.07~09.13"
$ws.Cells.Item(19, 3).Value = "20,000~23,000"
$ws.Cells.Item(19, 4).Value = "23000"
$ws.Cells.Item(19, 5).Value = 30000
$ws.Cells.Item(19, 6).Value = "미래에셋증권"

$ws.Cells.Item(20, 1).Value = "인스웨이브시스템즈"
$ws.Cells.Item(20, 2).Value = "2023.09.06~09.12"
$ws.Cells.Item(20, 3).Value = "20,000~24,000"
$ws.Cells.Item(20, 4).Value = "24000"
$ws.Cells.Item(20, 5).Value = 22000
$ws.Cells.Item(20, 6).Value = "신영증권"

$ws.Cells.Item(21, 1).Value = "아이엠티"
$ws.Cells.Item(21, 2).Value = "2023.09.06~09.12"
$ws.Cells.Item(21, 3).Value = "10,500~12,000"
$ws.Cells.Item(21, 4).Value = "14000"
$ws.Cells.Item(21, 5).Value = 16590
$ws.Cells.Item(21, 6).Value = "유안타증권,유진투자증권"
